$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 235
$ws1.Range("F5").Value = 2780
$ws1.Range("F6").Value = 1945
$ws1.Range("F7").Value = 375
$ws1.Range("F9").Value = 994
$ws1.Range("F11").Value = 29

# Sheet "全部类型" (all types)
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F4").Value = 235
$ws2.Range("F5").Value = 2780
$ws2.Range("F6").Value = 1945
$ws2.Range("F7").Value = 375
$ws2.Range("F10").Value = 994
$ws2.Range("F12").Value = 29
